$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the second student row (row 3, "Vo Duc Trong") ---
$ws.Rows.Item(3).Delete()

# --- Insert a new column before old column C for "Ngay Vao" (shifts old C..G -> D..H) ---
$ws.Columns.Item(3).Insert()

# --- Append two new columns (I, J) for "Thanh Tien" / "Trang Thai Dong Tien".
#     Copy header+data formatting from column H (last existing column) so the
#     new cells inherit the same borders / header styling. ---
$ws.Range("H1:H2").Copy()
$ws.Range("I1:I2").PasteSpecial(-4122)
$ws.Range("H1:H2").Copy()
$ws.Range("J1:J2").PasteSpecial(-4122)

# --- Header row (row 1) ---
$ws.Cells.Item(1,3).Value = "Ngày Vào"
$ws.Cells.Item(1,9).Value = "Thành Tiền"
$ws.Cells.Item(1,10).Value = "Trạng Thái Đóng Tiền"

# --- Data row (row 2): write every column explicitly so values line up with
#     their (corrected) headers. Force text format ("@") on cells that would
#     otherwise be auto-parsed as dates, then restore the normal bordered
#     cell style by copying it from a clean cell (A2). ---
$ws.Cells.Item(2,2).Value = "DH25PM"

$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = "2025-12-01"

$ws.Cells.Item(2,4).Value = "Lâm Huỳnh Phương"
$ws.Cells.Item(2,5).Value = "Nghiêm"
$ws.Cells.Item(2,6).Value = "Nữ"

$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "2006-10-23"

$ws.Cells.Item(2,9).NumberFormat = "@"
$ws.Cells.Item(2,9).Value = "2000000"

$ws.Cells.Item(2,10).Value = "Đã đóng"

# Re-apply the plain bordered data style (as used by A2) to every data cell in
# row 2 so none of them keep a stray "@"/date number format, and so any
# leftover style from the pre-edit layout (e.g. the old date format on what
# is now F2) is replaced with the normal bordered style.
$ws.Cells.Item(2,1).Copy()
$ws.Range("B2:J2").PasteSpecial(-4122)

# --- Column widths (A..J) ---
$widths = @(14,11,15,21,11,14,15,10,15,25)
for ($i = 1; $i -le 10; $i++) {
    $ws.Columns.Item($i).ColumnWidth = ($widths[$i-1] - 0.8333333333333333)
}

Write-Host "done"
